$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.727.07'
$ws.Range('E2').Value = '  +0.43%  '

$ws.Range('D3').Value = '1.848.00'
$ws.Range('E3').Value = '  +0.27%  '

$ws.Range('E4').Value = '  -0.07%  '

$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.73'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -0.40%  '

$ws.Range('E6').Value = '  +0.03%  '

$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4321'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  +1.23%  '

$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3658'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  +0.40%  '

$style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07339'
$ws.Range('D9').Style = $style
$ws.Range('E9').Value = '  +0.76%  '

$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8784'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  -1.73%  '

$ws.Range('E11').Value = '  +0.84%  '

$ws.Range('D12').Value = '1.861.19'
$ws.Range('E12').Value = '  +0.00%  '

$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.359'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  -0.40%  '

$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.540'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  -0.38%  '

$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06946'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  +1.05%  '

$ws.Range('E16').Value = '  -0.03%  '

$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.92'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  +2.06%  '

$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009013'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  +1.95%  '

$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  +0.02%  '

$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.37'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  -1.20%  '

$ws.Range('D21').Value = '27.669.20'
$ws.Range('E21').Value = '  +0.15%  '

$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.978'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  -0.06%  '

$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.33'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  -1.67%  '

$ws.Range('D24').Value = '2.098.41'
$ws.Range('E24').Value = '  +0.05%  '

$ws.Range('E25').Value = '  -2.53%  '

$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.85'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  +0.64%  '

$ws.Range('E27').Value = '  +1.58%  '

$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '120.05'
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  +6.27%  '

$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.254'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  +0.63%  '

$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.887'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  +2.65%  '

$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08906'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  +0.02%  '

$style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7600'
$ws.Range('D32').Style = $style
$ws.Range('E32').Value = '  -2.06%  '

$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.556'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  -0.45%  '

$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.964'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  -0.62%  '

$ws.Range('E35').Value = '  +2.65%  '

$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  +0.07%  '

$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.110'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  +0.81%  '

$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05445'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  +0.10%  '

$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01937'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  +0.61%  '

$ws.Range('E40').Value = '  +1.97%  '

$ws.Range('E41').Value = '  +0.68%  '

$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1665'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +1.29%  '

$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.611'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  -2.73%  '

$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.387'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  +1.89%  '

$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06556'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  -1.21%  '

$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.40'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +0.31%  '

$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '105.53'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -0.44%  '

$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4667'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  -0.90%  '

$ws.Range('E49').Value = '  +0.02%  '

$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.639'
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  +0.48%  '

$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.65'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  +0.31%  '
